$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Volume(1h) (E) columns store plain text values in this workbook
# (e.g. "1.00", "98.070.06", "  +0.45%  "), not numbers. Force those specific cells to
# Text format before writing so Excel does not auto-convert numeric-looking strings into
# real numbers (which would strip formatting like trailing zeros / thousands dots).
# (Applied as separate contiguous blocks since multi-area NumberFormat assignment only
# reliably touches the first area.)
$ws.Range("D2:D3").NumberFormat = "@"
$ws.Range("D5:D29").NumberFormat = "@"
$ws.Range("D31:D33").NumberFormat = "@"
$ws.Range("D35:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "98.070.06"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "3.294.21"
$ws.Range("E3").Value = "  -1.28%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "255.69"
$ws.Range("E5").Value = "  +3.79%  "
$ws.Range("D6").Value = "621.14"
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("D7").Value = "1.41"
$ws.Range("E7").Value = "  +25.02%  "
$ws.Range("D8").Value = "0.396"
$ws.Range("E8").Value = "  +1.96%  "
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("D10").Value = "0.892"
$ws.Range("E10").Value = "  +11.64%  "
$ws.Range("D11").Value = "3.288.69"
$ws.Range("E11").Value = "  -1.41%  "
$ws.Range("D12").Value = "0.197"
$ws.Range("E12").Value = "  -1.40%  "
$ws.Range("D13").Value = "37.93"
$ws.Range("E13").Value = "  +6.24%  "
$ws.Range("D14").Value = "97.764.65"
$ws.Range("E14").Value = "  +0.39%  "
$ws.Range("D15").Value = "0.0000245"
$ws.Range("E15").Value = "  -0.93%  "
$ws.Range("D16").Value = "3.906.00"
$ws.Range("E16").Value = "  -1.32%  "
$ws.Range("D17").Value = "5.49"
$ws.Range("E17").Value = "  -0.68%  "
$ws.Range("D18").Value = "3.307.86"
$ws.Range("E18").Value = "  -0.92%  "
$ws.Range("D19").Value = "3.50"
$ws.Range("E19").Value = "  -4.34%  "
$ws.Range("D20").Value = "15.06"
$ws.Range("E20").Value = "  -1.25%  "
$ws.Range("D21").Value = "476.91"
$ws.Range("E21").Value = "  -3.73%  "
$ws.Range("D22").Value = "6.04"
$ws.Range("E22").Value = "  +2.17%  "
$ws.Range("D23").Value = "0.0000203"
$ws.Range("E23").Value = "  -4.60%  "
$ws.Range("D24").Value = "9.27"
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").Value = "5.55"
$ws.Range("E25").Value = "  -2.30%  "
$ws.Range("D26").Value = "88.21"
$ws.Range("E26").Value = "  -0.68%  "
$ws.Range("D27").Value = "11.79"
$ws.Range("E27").Value = "  -2.52%  "
$ws.Range("D28").Value = "3.475.17"
$ws.Range("E28").Value = "  -0.83%  "
$ws.Range("D29").Value = "0.289"
$ws.Range("E29").Value = "  +19.64%  "
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("D31").Value = "0.186"
$ws.Range("E31").Value = "  +1.94%  "
$ws.Range("D32").Value = "0.131"
$ws.Range("E32").Value = "  +6.70%  "
$ws.Range("D33").Value = "9.74"
$ws.Range("E33").Value = "  +3.92%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").Value = "27.42"
$ws.Range("E35").Value = "  -1.05%  "
$ws.Range("D36").Value = "0.147"
$ws.Range("E36").Value = "  -5.99%  "
$ws.Range("D37").Value = "7.09"
$ws.Range("E37").Value = "  -6.11%  "
$ws.Range("D38").Value = "1.92"
$ws.Range("E38").Value = "  -1.18%  "
$ws.Range("B39").Value = "WhiteBITCoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D39").Value = "24.85"
$ws.Range("E39").Value = "  +0.49%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "490.49"
$ws.Range("E40").Value = "  -2.83%  "
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").Value = "0.454"
$ws.Range("E41").Value = "  +0.63%  "
$ws.Range("D42").Value = "3.67"
$ws.Range("E42").Value = "  +5.24%  "
$ws.Range("D43").Value = "1.23"
$ws.Range("E43").Value = "  -3.80%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").Value = "0.788"
$ws.Range("E44").Value = "  -0.64%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").Value = "3.19"
$ws.Range("E46").Value = "  -3.03%  "
$ws.Range("D47").Value = "159.45"
$ws.Range("E47").Value = "  -1.06%  "
$ws.Range("D48").Value = "1.89"
$ws.Range("E48").Value = "  -3.84%  "
$ws.Range("D49").Value = "0.828"
$ws.Range("E49").Value = "  +3.85%  "
$ws.Range("D50").Value = "4.57"
$ws.Range("E50").Value = "  -0.74%  "
$ws.Range("D51").Value = "45.64"
$ws.Range("E51").Value = "  +1.61%  "

Write-Output "Updated cryptos list"